$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.474.58"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.667.35"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'313.68"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").Value = "'0.3975"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("D8").Value = "'0.3917"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Value = "'52.07"
$ws.Range("E9").Value = "  +6.54%  "
$ws.Range("D10").Value = "'1.406"
$ws.Range("E10").Value = "  +3.48%  "
$ws.Range("D11").Value = "'0.9998"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "'0.08603"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "'24.42"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("E15").Value = "  +5.55%  "
$ws.Range("D16").Value = "'7.926"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("D17").Value = "1.663.65"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "'95.46"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "'0.06984"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'20.64"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'7.019"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "'0.9989"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'13.78"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "24.470.44"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "'2.426"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("E26").Value = "  +11.73%  "
$ws.Range("D27").Value = "'22.57"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'157.89"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'5.478"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").Value = "'142.95"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").Value = "'8.189"
$ws.Range("E31").Value = "  -8.90%  "
$ws.Range("D32").Value = "'2.520"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").Value = "1.847.54"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'1.073"
$ws.Range("E34").Value = "  +9.22%  "
$ws.Range("D35").Value = "'0.08292"
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("D36").Value = "'0.03044"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("D37").Value = "'6.917"
$ws.Range("E37").Value = "  -4.21%  "
$ws.Range("D38").Value = "'11.16"
$ws.Range("E38").Value = "  +11.32%  "
$ws.Range("D39").Value = "'0.2780"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("D40").Value = "'0.09254"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'13.92"
$ws.Range("E41").Value = "  +6.12%  "
$ws.Range("D42").Value = "'0.7761"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'1.446"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").Value = "'16.73"
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").Value = "'0.7144"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("D46").Value = "'2.547"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").Value = "'4.146"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("D48").Value = "'0.9994"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").Value = "'0.08468"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "'136.83"
$ws.Range("E50").Value = "  +2.00%  "
$ws.Range("D51").Value = "'1.274"
$ws.Range("E51").Value = "  +0.84%  "
